$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.863.72"
$ws.Range("E2").Value = "  -2.18%  "
$ws.Range("D3").Value = "3.588.43"
$ws.Range("E3").Value = "  -2.05%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.39"
$ws.Range("E5").Value = "  -1.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.10"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("E7").Value = "  -3.87%  "
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.662"
$ws.Range("E9").Value = "  -7.18%  "
$ws.Range("E10").Value = "  -11.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.56"
$ws.Range("E11").Value = "  -6.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000248"
$ws.Range("E12").Value = "  -15.02%  "
$ws.Range("E13").Value = "  -7.68%  "
$ws.Range("D14").Value = "4.168.21"
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("D15").Value = "3.591.18"
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("D17").Value = "66.575.85"
$ws.Range("E17").Value = "  -2.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.15"
$ws.Range("E18").Value = "  -5.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.06"
$ws.Range("E19").Value = "  -5.70%  "
$ws.Range("E20").Value = "  -6.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "387.50"
$ws.Range("E21").Value = "  -5.34%  "
$ws.Range("E22").Value = "  -7.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "83.88"
$ws.Range("E23").Value = "  -5.22%  "
$ws.Range("E24").Value = "  -6.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.04"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.03"
$ws.Range("E26").Value = "  -5.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.09"
$ws.Range("E27").Value = "  -7.15%  "
$ws.Range("E28").Value = "  -9.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.82"
$ws.Range("E29").Value = "  -6.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.75"
$ws.Range("E30").Value = "  -5.67%  "
$ws.Range("E31").Value = "  -8.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "65.36"
$ws.Range("E32").Value = "  +1.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.72"
$ws.Range("E33").Value = "  -5.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "588.48"
$ws.Range("E34").Value = "  -2.44%  "
$ws.Range("E35").Value = "  -5.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "40.63"
$ws.Range("E36").Value = "  -5.99%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E39").Value = "  -7.71%  "
$ws.Range("D40").Value = "0.0₃0733"
$ws.Range("E40").Value = "  -17.80%  "
$ws.Range("E41").Value = "  -4.10%  "
$ws.Range("E42").Value = "  -9.86%  "
$ws.Range("E43").Value = "  -7.42%  "
$ws.Range("D44").Value = "2.659.94"
$ws.Range("E44").Value = "  -1.98%  "
$ws.Range("E45").Value = "  -13.35%  "
$ws.Range("E46").Value = "  -4.80%  "
$ws.Range("E47").Value = "  -3.65%  "
$ws.Range("E48").Value = "  -7.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "134.18"
$ws.Range("E49").Value = "  -5.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.12"
$ws.Range("E50").Value = "  -9.97%  "
$ws.Range("E51").Value = "  -8.05%  "
